$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value would otherwise be auto-converted to a number by Excel;
# force them to keep a text format so the stored value matches the literal string.
$textCells = @("D5","D6","D10","D11","D16","D19","D20","D21","D22","D23","D24","D26","D29","D32","D33","D35","D37","D39","D40","D42","D43","D45","D47","D49","D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "55.841.45"
$ws.Range("E2").Value = "  +11.08%  "
$ws.Range("D3").Value = "2.523.80"
$ws.Range("E3").Value = "  +15.35%  "
$ws.Range("E4").Value = "  -0.12%  "
$ws.Range("D5").Value = "491.38"
$ws.Range("E5").Value = "  +19.43%  "
$ws.Range("D6").Value = "142.43"
$ws.Range("E6").Value = "  +27.58%  "
$ws.Range("E7").Value = "  -0.20%  "
$ws.Range("E8").Value = "  +16.99%  "
$ws.Range("D9").Value = "2.516.61"
$ws.Range("E9").Value = "  +16.72%  "
$ws.Range("D10").Value = "0.0996"
$ws.Range("E10").Value = "  +20.07%  "
$ws.Range("D11").Value = "5.55"
$ws.Range("E11").Value = "  +13.25%  "
$ws.Range("E12").Value = "  +19.17%  "
$ws.Range("E13").Value = "  +4.84%  "
$ws.Range("D14").Value = "2.965.86"
$ws.Range("E14").Value = "  +16.23%  "
$ws.Range("D15").Value = "55.878.03"
$ws.Range("E15").Value = "  +11.12%  "
$ws.Range("D16").Value = "20.91"
$ws.Range("E16").Value = "  +18.68%  "
$ws.Range("E17").Value = "  +27.08%  "
$ws.Range("D18").Value = "2.526.83"
$ws.Range("E18").Value = "  +13.82%  "
$ws.Range("D19").Value = "4.44"
$ws.Range("E19").Value = "  +21.34%  "
$ws.Range("D20").Value = "325.27"
$ws.Range("E20").Value = "  +16.71%  "
$ws.Range("D21").Value = "10.10"
$ws.Range("E21").Value = "  +23.25%  "
$ws.Range("D22").Value = "0.999"
$ws.Range("E22").Value = "  +0.10%  "
$ws.Range("D23").Value = "5.80"
$ws.Range("E23").Value = "  +18.18%  "
$ws.Range("D24").Value = "58.64"
$ws.Range("E24").Value = "  +14.61%  "
$ws.Range("E25").Value = "  +31.25%  "
$ws.Range("D26").Value = "0.413"
$ws.Range("E26").Value = "  +21.81%  "
$ws.Range("E27").Value = "  -0.01%  "
$ws.Range("D28").Value = "2.629.89"
$ws.Range("E28").Value = "  +14.55%  "
$ws.Range("D29").Value = "7.54"
$ws.Range("E29").Value = "  +17.68%  "
$ws.Range("D30").Value = "0.0₃0810"
$ws.Range("E30").Value = "  +31.74%  "
$ws.Range("E31").Value = "  +0.04%  "
$ws.Range("D32").Value = "150.73"
$ws.Range("E32").Value = "  +9.18%  "
$ws.Range("D33").Value = "18.36"
$ws.Range("E33").Value = "  +15.30%  "
$ws.Range("E34").Value = "  +22.57%  "
$ws.Range("D35").Value = "5.25"
$ws.Range("E35").Value = "  +19.88%  "
$ws.Range("E36").Value = "  +19.39%  "
$ws.Range("D37").Value = "3.76"
$ws.Range("E37").Value = "  +17.47%  "
$ws.Range("E38").Value = "  +22.41%  "
$ws.Range("D39").Value = "34.38"
$ws.Range("E39").Value = "  +12.96%  "
$ws.Range("D40").Value = "0.618"
$ws.Range("E40").Value = "  +27.99%  "
$ws.Range("E41").Value = "  +22.11%  "
$ws.Range("B42").Value = "FirstDigitalUSD"
$ws.Range("C42").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D42").Value = "0.999"
$ws.Range("E42").Value = "  +0.37%  "
$ws.Range("B43").Value = "Filecoin"
$ws.Range("C43").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D43").Value = "3.46"
$ws.Range("E43").Value = "  +19.06%  "
$ws.Range("E44").Value = "  +18.81%  "
$ws.Range("B45").Value = "RenderToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D45").Value = "4.79"
$ws.Range("E45").Value = "  +34.26%  "
$ws.Range("B46").Value = "Maker"
$ws.Range("C46").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D46").Value = "2.019.12"
$ws.Range("E46").Value = "  +13.46%  "
$ws.Range("D47").Value = "259.53"
$ws.Range("E47").Value = "  +55.56%  "
$ws.Range("E48").Value = "  +17.99%  "
$ws.Range("D49").Value = "10.15"
$ws.Range("E49").Value = "  -0.31%  "
$ws.Range("E50").Value = "  +18.70%  "
$ws.Range("D51").Value = "17.81"
$ws.Range("E51").Value = "  +22.85%  "
